# Add package information to Telegram base classes.
# Inserts two new rows (14 and 15) into the "valueObject" sheet's
# inheritance block, directly below the inherited class name row (13),
# carrying the 名前空間 (namespace) and パッケージ (package) of the base
# class ApiTelegram. Everything from the former row 14 onward shifts down
# by two rows (old 14 -> 16, ... old 32 -> 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("valueObject")

# Insert two blank rows right after row 13 (the inherited class-name row),
# pushing the old rows 14-32 down to 16-34.
$ws.Rows("14:15").Insert()

# --- Formatting: copy the look of analogous "common section" rows -------
# Row 14 (namespace) mirrors the common-section namespace row (row 7).
$ws.Range("A7:D7").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)
# Column E on row 14 uses the same fill/top-border treatment as D7.
$ws.Range("D7").Copy()
$ws.Range("E14").PasteSpecial(-4122)
# Column F matches the class-name row directly above (row 13).
$ws.Range("F13").Copy()
$ws.Range("F14").PasteSpecial(-4122)

# Row 15 (package) mirrors the common-section package/class rows.
$ws.Range("A7:B7").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$ws.Range("C6:D6").Copy()
$ws.Range("C15:D15").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F13").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Values ---------------------------------------------------------------
$ws.Range("A14").Value = "名前空間"
$ws.Range("C14").Value = "\blanco\sample\valueobject"

$ws.Range("A15").Value = "パッケージ"
$ws.Range("C15").Value = "blanco.rest.valueobject"
